$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5: dilation (G5) changes from 1 to 2 (B5/C5 recalc automatically since
# they are formulas driven by shared-formula groups referencing this row).
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = 2

# ---------------------------------------------------------------------------
# Rows 7-12: the layer pattern (ConvTrans / Upsampling alternating) is
# shifted by one row and a new final ConvTrans row is appended. We rewrite
# every cell in A7:H12 explicitly to land on the target content.
#
# First, repaint the "disabled" shaded fill (style used by the Upsampling
# rows, which do not use Padding/Stride/Kernel/Dilation/Output Padding) onto
# its new rows (7, 9, 11), and clear that shading from the rows that become
# ConvTrans rows (8, 10, 12), using the *original* layout (row 8 = shaded,
# row 7 = plain) as the format source before any values are changed.
# ---------------------------------------------------------------------------
$ws.Range("D8:H8").Copy() | Out-Null
$ws.Range("D9:H9").PasteSpecial(-4122) | Out-Null
$ws.Range("D8:H8").Copy() | Out-Null
$ws.Range("D11:H11").PasteSpecial(-4122) | Out-Null

$ws.Range("D7:H7").Copy() | Out-Null
$ws.Range("D10:H10").PasteSpecial(-4122) | Out-Null
$ws.Range("D7:H7").Copy() | Out-Null
$ws.Range("D12:H12").PasteSpecial(-4122) | Out-Null

$ws.Range("D8:H8").Copy() | Out-Null
$ws.Range("D7:H7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7:H7").Copy() | Out-Null
$ws.Range("D8:H8").PasteSpecial(-4122) | Out-Null

# Row 7 -> Upsampling (was ConvTrans)
$ws.Range("A7").Value = "Upsampling"
$ws.Range("B7").Formula = "=C7*2"
$ws.Range("C7").Formula = "=B6"
$ws.Range("D7:H7").ClearContents()

# Row 8 -> ConvTrans (was Upsampling)
$ws.Range("A8").Value = "ConvTrans"
$ws.Range("B8").Formula = "=(C8-1)*E8-2*D8+G8*(F8-1)+H8+1"
$ws.Range("C8").Formula = "=B7"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0

# Row 9 -> Upsampling (was ConvTrans)
$ws.Range("A9").Value = "Upsampling"
$ws.Range("B9").Formula = "=C9*2"
$ws.Range("C9").Formula = "=B8"
$ws.Range("D9:H9").ClearContents()

# Row 10 -> ConvTrans (was Upsampling)
$ws.Range("A10").Value = "ConvTrans"
$ws.Range("B10").Formula = "=(C10-1)*E10-2*D10+G10*(F10-1)+H10+1"
$ws.Range("C10").Formula = "=B9"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0

# Row 11 -> Upsampling (was ConvTrans)
$ws.Range("A11").Value = "Upsampling"
$ws.Range("B11").Formula = "=C11*2"
$ws.Range("C11").Formula = "=B10"
$ws.Range("D11:H11").ClearContents()

# Row 12 -> ConvTrans (was Upsampling) - also gets a new Output Padding of 1
$ws.Range("A12").Value = "ConvTrans"
$ws.Range("B12").Formula = "=(C12-1)*E12-2*D12+G12*(F12-1)+H12+1"
$ws.Range("C12").Formula = "=B11"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

# ---------------------------------------------------------------------------
# New row 16: a validation warning message in column F.
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = "output padding must be smaller then stride or dilation"

# ---------------------------------------------------------------------------
# Cosmetic sheet-level updates: selection, column width.
# ---------------------------------------------------------------------------
$ws.Range("H13").Select()
$ws.Columns("B").ColumnWidth = 21
